$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.426.59"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.563.50"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "286.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3342"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07417"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.935"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.888"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "1.562.08"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001106"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06701"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.337"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "22.406.19"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.551"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.007"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("D31").Value = "1.736.70"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.062"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.132"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.989"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.607"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08237"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02391"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06386"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.300"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2212"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.348"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6074"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5746"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.212"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.46%  "
